# Append 12 more rows (A=204..215, matching B values) to the existing
# "특징값" (normalized feature value) series that runs from row 2..205.
# Column A keeps the bold/bordered/centered style (s="1") used by the
# existing index column; column B keeps the default style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
  @(204, "1.162264728904461E-16"),
  @(205, "2.087976256539357E-16"),
  @(206, "9.853229343548264E-17"),
  @(207, "-3.700743415417188E-17"),
  @(208, "-3.469446951953614E-18"),
  @(209, "-2.577303450022685E-17"),
  @(210, "-1.225871256356944E-16"),
  @(211, "-1.554312234475219E-16"),
  @(212, "4.510281037539698E-17"),
  @(213, "-6.47630097698008E-17"),
  @(214, "-9.71445146547012E-17"),
  @(215, "0")
)

$row = 206
foreach ($pair in $newRows) {
  $ws.Cells.Item($row, 1).Value = $pair[0]
  $ws.Cells.Item($row, 2).Value = [double]$pair[1]
  $row = $row + 1
}

# Copy column A's existing style (bold, bordered, centered) down onto the
# newly added A cells so they match the rest of the index column.
$ws.Range("A205").Copy()
$ws.Range("A206:A217").PasteSpecial(-4122)
